$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (30, 31) continuing the existing daily series.
# Use the raw date serial numbers (matching column A's existing values) and
# copy row 29's date formatting so the new cells share its style.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(30, 1).Value = 44785
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0

$ws.Cells.Item(31, 1).Value = 44786
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0

# Match the author's final selection state.
$ws.Range("H33").Select()
